$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (it currently sits between "products "
#    and "and to innovate with our live streaming."). Word will later decide
#    where the new "last edit" location is, and we re-create the bookmark
#    there ourselves (see below), matching the target document exactly.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Rewrite "... Our platform will suggest a minimal way of messaging." as
#    "... Our platform will suggest a basic messaging system." split across
#    several runs (mirroring the target OOXML), with a new "_GoBack" bookmark
#    placed immediately before the final "." run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "There are lots of text-messaging systems, such as WhatsApp, Telegram, etc. Our platform will suggest a minimal way of messaging.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

$oldStart = $rng.Start
$oldEnd = $rng.End

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document ' + $ns + '><w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">There are lots of text-messaging systems, such as WhatsApp, Telegram, etc. Our platform will suggest a </w:t></w:r>' + `
  '<w:bookmarkStart w:id="101" w:name="zzTmpSplit1"/><w:bookmarkEnd w:id="101"/>' + `
  '<w:r><w:t xml:space="preserve">basic </w:t></w:r>' + `
  '<w:bookmarkStart w:id="102" w:name="zzTmpSplit2"/><w:bookmarkEnd w:id="102"/>' + `
  '<w:r><w:t>messaging</w:t></w:r>' + `
  '<w:bookmarkStart w:id="103" w:name="zzTmpSplit3"/><w:bookmarkEnd w:id="103"/>' + `
  '<w:r><w:t xml:space="preserve"> system</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML places the new, already-split runs right after the matched
# (old) text, inside the same paragraph (preserving its pPr/identity).
$rng.InsertXML($xml)

# Remove the old, now-duplicated sentence text.
$oldRange = $d.Range($oldStart, $oldEnd)
$oldRange.Text = ""

# Drop the temporary bookmarks that kept our new runs from being silently
# re-merged back together; only "_GoBack" should remain in the document.
foreach ($name in @("zzTmpSplit1", "zzTmpSplit2", "zzTmpSplit3")) {
    if ($d.Bookmarks.Exists($name)) {
        $d.Bookmarks($name).Delete()
    }
}
